$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the data set.
# Delete the higher-numbered row first (SC 92, row 28) so the lower
# row index (RM 232, row 26) is unaffected by the shift.
$ws.Rows.Item(28).Delete() | Out-Null
$ws.Rows.Item(26).Delete() | Out-Null

# Apply the individual cell value changes (some values filled in,
# some cleared back to blank/missing) on the resulting 33-row table.
$ws.Range("F2").Value = ""
$ws.Range("C3").Value = 11.2
$ws.Range("F3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = 17.66
$ws.Range("E6").Value = -5.7
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("E12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = -5.4
$ws.Range("E19").Value = -6.5
$ws.Range("E20").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("E27").Value = -10
$ws.Range("E28").Value = -5.9
$ws.Range("F31").Value = 17.18
$ws.Range("C32").Value = 10.5
